$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 156.656447
$ws.Range("H2").Value = 469.969341
$ws.Range("I2").Value = 0.0671576211124673
$ws.Range("J2").Value = 0.0671576211124673
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.664391
$ws.Range("N2").Value = 4.993173000000001
$ws.Range("O2").Value = 0.3990511495040125
$ws.Range("P2").Value = 0.3990511495040125
$ws.Range("Q2").Value = 260.737580478777
$ws.Range("R2").Value = 2346.638224308993
$ws.Range("S2").Value = 0.02679932590288502
$ws.Range("T2").Value = 0.02679932590288502

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 156.656447
$ws.Range("H3").Value = 469.969341
$ws.Range("I3").Value = 0.0671576211124673
$ws.Range("J3").Value = 0.0671576211124673
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.437958
$ws.Range("N3").Value = 1.313874
$ws.Range("O3").Value = 0.1050039584054939
$ws.Range("P3").Value = 0.1050039584054938
$ws.Range("Q3").Value = 68.60894421522599
$ws.Range("R3").Value = 617.480497937034
$ws.Range("S3").Value = 0.007051816053905432
$ws.Range("T3").Value = 0.007051816053905431

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 156.656447
$ws.Range("H4").Value = 469.969341
$ws.Range("I4").Value = 0.0671576211124673
$ws.Range("J4").Value = 0.0671576211124673
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.439215333333333
$ws.Range("N4").Value = 4.317646
$ws.Range("O4").Value = 0.3450634695516061
$ws.Range("P4").Value = 0.3450634695516061
$ws.Range("Q4").Value = 225.4623605879206
$ws.Range("R4").Value = 2029.161245291286
$ws.Range("S4").Value = 0.02317364174790016
$ws.Range("T4").Value = 0.02317364174790016

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 156.656447
$ws.Range("H5").Value = 469.969341
$ws.Range("I5").Value = 0.0671576211124673
$ws.Range("J5").Value = 0.0671576211124673
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.629307
$ws.Range("N5").Value = 1.887921
$ws.Range("O5").Value = 0.1508814225388875
$ws.Range("P5").Value = 0.1508814225388875
$ws.Range("Q5").Value = 98.58499869222898
$ws.Range("R5").Value = 887.2649882300609
$ws.Range("S5").Value = 0.01013283740777669
$ws.Range("T5").Value = 0.01013283740777669

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 42.300692
$ws.Range("H6").Value = 126.902076
$ws.Range("I6").Value = 0.01813403725498241
$ws.Range("J6").Value = 0.01813403725498241
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.664391
$ws.Range("N6").Value = 4.993173000000001
$ws.Range("O6").Value = 0.3990511495040125
$ws.Range("P6").Value = 0.3990511495040125
$ws.Range("Q6").Value = 70.404891058572
$ws.Range("R6").Value = 633.6440195271481
$ws.Range("S6").Value = 0.007236408411749318
$ws.Range("T6").Value = 0.007236408411749318

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 42.300692
$ws.Range("H7").Value = 126.902076
$ws.Range("I7").Value = 0.01813403725498241
$ws.Range("J7").Value = 0.01813403725498241
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.437958
$ws.Range("N7").Value = 1.313874
$ws.Range("O7").Value = 0.1050039584054939
$ws.Range("P7").Value = 0.1050039584054938
$ws.Range("Q7").Value = 18.525926466936
$ws.Range("R7").Value = 166.733338202424
$ws.Range("S7").Value = 0.001904145693645849
$ws.Range("T7").Value = 0.001904145693645848

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 42.300692
$ws.Range("H8").Value = 126.902076
$ws.Range("I8").Value = 0.01813403725498241
$ws.Range("J8").Value = 0.01813403725498241
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.439215333333333
$ws.Range("N8").Value = 4.317646
$ws.Range("O8").Value = 0.3450634695516061
$ws.Range("P8").Value = 0.3450634695516061
$ws.Range("Q8").Value = 60.87980453701066
$ws.Range("R8").Value = 547.918240833096
$ws.Range("S8").Value = 0.006257393812182311
$ws.Range("T8").Value = 0.006257393812182311

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 42.300692
$ws.Range("H9").Value = 126.902076
$ws.Range("I9").Value = 0.01813403725498241
$ws.Range("J9").Value = 0.01813403725498241
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.629307
$ws.Range("N9").Value = 1.887921
$ws.Range("O9").Value = 0.1508814225388875
$ws.Range("P9").Value = 0.1508814225388875
$ws.Range("Q9").Value = 26.620121580444
$ws.Range("R9").Value = 239.581094223996
$ws.Range("S9").Value = 0.002736089337404928
$ws.Range("T9").Value = 0.002736089337404928

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2110.189616
$ws.Range("H10").Value = 6330.568848
$ws.Range("I10").Value = 0.9046248489651427
$ws.Range("J10").Value = 0.9046248489651426
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.664391
$ws.Range("N10").Value = 4.993173000000001
$ws.Range("O10").Value = 0.3990511495040125
$ws.Range("P10").Value = 0.3990511495040125
$ws.Range("Q10").Value = 3512.180605163856
$ws.Range("R10").Value = 31609.62544647471
$ws.Range("S10").Value = 0.3609915858494339
$ws.Range("T10").Value = 0.3609915858494339

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2110.189616
$ws.Range("H11").Value = 6330.568848
$ws.Range("I11").Value = 0.9046248489651427
$ws.Range("J11").Value = 0.9046248489651426
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.437958
$ws.Range("N11").Value = 1.313874
$ws.Range("O11").Value = 0.1050039584054939
$ws.Range("P11").Value = 0.1050039584054938
$ws.Range("Q11").Value = 924.1744238441281
$ws.Range("R11").Value = 8317.569814597151
$ws.Range("S11").Value = 0.094989190013312
$ws.Range("T11").Value = 0.09498919001331198

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2110.189616
$ws.Range("H12").Value = 6330.568848
$ws.Range("I12").Value = 0.9046248489651427
$ws.Range("J12").Value = 0.9046248489651426
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.439215333333333
$ws.Range("N12").Value = 4.317646
$ws.Range("O12").Value = 0.3450634695516061
$ws.Range("P12").Value = 0.3450634695516061
$ws.Range("Q12").Value = 3037.017251587979
$ws.Range("R12").Value = 27333.15526429181
$ws.Range("S12").Value = 0.3121529890265097
$ws.Range("T12").Value = 0.3121529890265097

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2110.189616
$ws.Range("H13").Value = 6330.568848
$ws.Range("I13").Value = 0.9046248489651427
$ws.Range("J13").Value = 0.9046248489651426
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.629307
$ws.Range("N13").Value = 1.887921
$ws.Range("O13").Value = 0.1508814225388875
$ws.Range("P13").Value = 0.1508814225388875
$ws.Range("Q13").Value = 1327.957096676112
$ws.Range("R13").Value = 11951.61387008501
$ws.Range("S13").Value = 0.136491084075887
$ws.Range("T13").Value = 0.136491084075887

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 23.52144266666667
$ws.Range("H14").Value = 70.564328
$ws.Range("I14").Value = 0.01008349266740757
$ws.Range("J14").Value = 0.01008349266740757
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.664391
$ws.Range("N14").Value = 4.993173000000001
$ws.Range("O14").Value = 0.3990511495040125
$ws.Range("P14").Value = 0.3990511495040125
$ws.Range("Q14").Value = 39.14887748141601
$ws.Range("R14").Value = 352.339897332744
$ws.Range("S14").Value = 0.004023829339944273
$ws.Range("T14").Value = 0.004023829339944272

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 23.52144266666667
$ws.Range("H15").Value = 70.564328
$ws.Range("I15").Value = 0.01008349266740757
$ws.Range("J15").Value = 0.01008349266740757
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.437958
$ws.Range("N15").Value = 1.313874
$ws.Range("O15").Value = 0.1050039584054939
$ws.Range("P15").Value = 0.1050039584054938
$ws.Range("Q15").Value = 10.301403987408
$ws.Range("R15").Value = 92.712635886672
$ws.Range("S15").Value = 0.001058806644630567
$ws.Range("T15").Value = 0.001058806644630567

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 23.52144266666667
$ws.Range("H16").Value = 70.564328
$ws.Range("I16").Value = 0.01008349266740757
$ws.Range("J16").Value = 0.01008349266740757
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.439215333333333
$ws.Range("N16").Value = 4.317646
$ws.Range("O16").Value = 0.3450634695516061
$ws.Range("P16").Value = 0.3450634695516061
$ws.Range("Q16").Value = 33.85242094798756
$ws.Range("R16").Value = 304.671788531888
$ws.Range("S16").Value = 0.003479444965013835
$ws.Range("T16").Value = 0.003479444965013835

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 23.52144266666667
$ws.Range("H17").Value = 70.564328
$ws.Range("I17").Value = 0.01008349266740757
$ws.Range("J17").Value = 0.01008349266740757
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.629307
$ws.Range("N17").Value = 1.887921
$ws.Range("O17").Value = 0.1508814225388875
$ws.Range("P17").Value = 0.1508814225388875
$ws.Range("Q17").Value = 14.802208520232
$ws.Range("R17").Value = 133.219876682088
$ws.Range("S17").Value = 0.001521411717818896
$ws.Range("T17").Value = 0.001521411717818896
